$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Means ----
$ws1 = $wb.Worksheets.Item("Means")

# New headers for new columns
$ws1.Range("F1").Value = "Within 5 miles of HFC production facility"
$ws1.Range("G1").Value = "Within 10 miles of HFC production facility"

# New column values for existing rows 2-8 (unchanged rows get new F/G values)
$ws1.Range("F2").Value = 62
$ws1.Range("G2").Value = 65

$ws1.Range("F3").Value = 36
$ws1.Range("G3").Value = 28

$ws1.Range("F4").Value = 3
$ws1.Range("G4").Value = 7.2

$ws1.Range("F5").Value = 2.9
$ws1.Range("G5").Value = 5.1

$ws1.Range("F6").Value = 80
$ws1.Range("G6").Value = 79

$ws1.Range("F7").Value = 2.8
$ws1.Range("G7").Value = 5.7

$ws1.Range("F8").Value = 5.5
$ws1.Range("G8").Value = 5

# Row 9 (Total Cancer Risk) values changed + new columns
$ws1.Range("B9").Value = 26
$ws1.Range("C9").Value = 39
$ws1.Range("D9").Value = 110
$ws1.Range("E9").Value = 120
$ws1.Range("F9").Value = 120
$ws1.Range("G9").Value = 79

# Row 10 (Total Respiratory) values changed + new columns
$ws1.Range("B10").Value = 0.32
$ws1.Range("C10").Value = 0.43
$ws1.Range("D10").Value = 0.52
$ws1.Range("E10").Value = 0.53
$ws1.Range("F10").Value = 0.53
$ws1.Range("G10").Value = 0.51

# ---- Sheet 2: Standard Deviations ----
$ws2 = $wb.Worksheets.Item("Standard Deviations")

$ws2.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$ws2.Range("G1").Value = "Within 10 mile of HFC production facility SD"

$ws2.Range("F2").Value = 26
$ws2.Range("G2").Value = 26

$ws2.Range("F3").Value = 26
$ws2.Range("G3").Value = 26

$ws2.Range("F4").Value = 2.7
$ws2.Range("G4").Value = 6.7

$ws2.Range("F5").Value = 8.5
$ws2.Range("G5").Value = 6.4

$ws2.Range("F6").Value = 33
$ws2.Range("G6").Value = 32

$ws2.Range("F7").Value = 2.3
$ws2.Range("G7").Value = 7.3

$ws2.Range("F8").Value = 9.6
$ws2.Range("G8").Value = 8

$ws2.Range("B9").Value = 8.6
$ws2.Range("C9").Value = 24
$ws2.Range("D9").Value = 64
$ws2.Range("E9").Value = 61
$ws2.Range("F9").Value = 57
$ws2.Range("G9").Value = 32

$ws2.Range("B10").Value = 0.14
$ws2.Range("C10").Value = 0.084
$ws2.Range("D10").Value = 0.05
$ws2.Range("E10").Value = 0.052
$ws2.Range("F10").Value = 0.085
$ws2.Range("G10").Value = 0.066
